# Add a "Skill Description" column (full skill names) between SkillCode (A)
# and SFIA Level (old B) columns, per commit:
#   "skill full names added to PDP outputs and spreadsheets"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing SFIA Level / Keycode / Description columns one place right
# and make room for the new "Skill Description" column at B.
$ws.Columns.Item(2).Insert()

# Header
$ws.Range("B1").Value = "Skill Description"

# Autonomy block
$ws.Range("B2").Value = "Autonomy"
$ws.Range("B3").Value = "Autonomy"

# Influence block
$ws.Range("B4").Value = "Influence"
$ws.Range("B5").Value = "Influence"
$ws.Range("B6").Value = "Influence"
$ws.Range("B7").Value = "Influence"

# Complexity block
$ws.Range("B8").Value = "Complexity"
$ws.Range("B9").Value = "Complexity"
$ws.Range("B10").Value = "Complexity"

# Knowledge block
$ws.Range("B11").Value = "Knowledge"
$ws.Range("B12").Value = "Knowledge"
$ws.Range("B13").Value = "Knowledge"

# Row 14 stays blank (separator row)

# SUPP -> Supplier management
$ws.Range("B15").Value = "Supplier management"
$ws.Range("B16").Value = "Supplier management"
$ws.Range("B17").Value = "Supplier management"
$ws.Range("B18").Value = "Supplier management"
$ws.Range("B19").Value = "Supplier management"

# RLMT -> Stakeholder relationship management
$ws.Range("B20").Value = "Stakeholder relationship management"
$ws.Range("B21").Value = "Stakeholder relationship management"
$ws.Range("B22").Value = "Stakeholder relationship management"
$ws.Range("B23").Value = "Stakeholder relationship management"
